$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix production year for row 198 (2020 -> 2019)
$ws.Range("C198").Value = 2019

# Populate technical specification columns (E:T) for rows 167-208.
# Rows 167-201 receive real data; rows 202-208 are touched (span grows)
# but carry no values, matching the source edit.
$arr = New-Object 'object[,]' 42,16
$arr[0,0] = 2.5
$arr[0,1] = 4
$arr[0,2] = 300
$arr[0,3] = 5
$arr[0,4] = 4
$arr[0,5] = 5.2
$arr[0,6] = 255
$arr[0,7] = 407
$arr[0,8] = 1500
$arr[0,9] = 460
$arr[0,10] = 60
$arr[0,11] = 460
$arr[0,12] = 180
$arr[0,13] = 148
$arr[0,14] = 265
$arr[0,15] = 6
$arr[1,0] = 4
$arr[1,1] = 8
$arr[1,2] = 600
$arr[1,3] = 2
$arr[1,4] = 2
$arr[1,5] = 3.7
$arr[1,6] = 304
$arr[1,7] = 750
$arr[1,8] = 2195
$arr[1,9] = 0
$arr[1,10] = 90
$arr[1,11] = 480
$arr[1,12] = 194
$arr[1,13] = 139
$arr[1,14] = 275
$arr[1,15] = 8
$arr[2,0] = 3.5
$arr[2,1] = 6
$arr[2,2] = 350
$arr[2,3] = 2
$arr[2,4] = 2
$arr[2,5] = 4
$arr[2,6] = 274
$arr[2,7] = 400
$arr[2,8] = 1161
$arr[2,9] = 134
$arr[2,10] = 40
$arr[2,11] = 308
$arr[2,12] = 180
$arr[2,13] = 113
$arr[2,14] = 237
$arr[2,15] = 6
$arr[3,0] = 5.4
$arr[3,1] = 8
$arr[3,2] = 305
$arr[3,3] = 2
$arr[3,4] = 2
$arr[3,5] = 8.2
$arr[3,6] = 237
$arr[3,7] = 624
$arr[3,8] = 1525
$arr[3,9] = 297
$arr[3,10] = 76
$arr[3,11] = 445
$arr[3,12] = 177
$arr[3,13] = 127
$arr[3,14] = 143
$arr[3,15] = 4
$arr[4,0] = 3.8
$arr[4,1] = 6
$arr[4,2] = 306
$arr[4,3] = 4
$arr[4,4] = 2
$arr[4,5] = 6
$arr[4,6] = 240
$arr[4,7] = 361
$arr[4,8] = 1560
$arr[4,9] = 332
$arr[4,10] = 65
$arr[4,11] = 463
$arr[4,12] = 186
$arr[4,13] = 138
$arr[4,14] = 282
$arr[4,15] = 6
$arr[5,0] = 4.2
$arr[5,1] = 8
$arr[5,2] = 350
$arr[5,3] = 5
$arr[5,4] = 4
$arr[5,5] = 5.9
$arr[5,6] = 250
$arr[5,7] = 440
$arr[5,8] = 1770
$arr[5,9] = 546
$arr[5,10] = 80
$arr[5,11] = 492
$arr[5,12] = 186
$arr[5,13] = 146
$arr[5,14] = 284
$arr[5,15] = 6
$arr[6,0] = 0
$arr[6,1] = 2
$arr[6,2] = 646
$arr[6,3] = 4
$arr[6,4] = 4
$arr[6,5] = 3.3
$arr[6,6] = 250
$arr[6,7] = 830
$arr[6,8] = 2347
$arr[6,9] = 350
$arr[6,10] = 0
$arr[6,11] = 499
$arr[6,12] = 196
$arr[6,13] = 141
$arr[6,14] = 290
$arr[6,15] = 2
$arr[7,0] = 4
$arr[7,1] = 8
$arr[7,2] = 490
$arr[7,3] = 2
$arr[7,4] = 2
$arr[7,5] = 3.9
$arr[7,6] = 249
$arr[7,7] = 500
$arr[7,8] = 1113
$arr[7,9] = 0
$arr[7,10] = 120
$arr[7,11] = 501
$arr[7,12] = 195
$arr[7,13] = 125
$arr[7,14] = 275
$arr[7,15] = 6
$arr[8,0] = 6
$arr[8,1] = 12
$arr[8,2] = 592
$arr[8,3] = 2
$arr[8,4] = 2
$arr[8,5] = 3.8
$arr[8,6] = 330
$arr[8,7] = 700
$arr[8,8] = 1000
$arr[8,9] = 0
$arr[8,10] = 90
$arr[8,11] = 486
$arr[8,12] = 195
$arr[8,13] = 110
$arr[8,14] = 267
$arr[8,15] = 6
$arr[9,0] = 3
$arr[9,1] = 6
$arr[9,2] = 211
$arr[9,3] = 5
$arr[9,4] = 2
$arr[9,5] = 8.4
$arr[9,6] = 243
$arr[9,7] = 290
$arr[9,8] = 1687
$arr[9,9] = 400
$arr[9,10] = 67
$arr[9,11] = 482
$arr[9,12] = 187
$arr[9,13] = 140
$arr[9,14] = 273
$arr[9,15] = 6
$arr[10,0] = 7.5
$arr[10,1] = 8
$arr[10,2] = 305
$arr[10,3] = 2
$arr[10,4] = 2
$arr[10,5] = 6.1
$arr[10,6] = 176
$arr[10,7] = 563
$arr[10,8] = 1694
$arr[10,9] = 205
$arr[10,10] = 64
$arr[10,11] = 500
$arr[10,12] = 185
$arr[10,13] = 128
$arr[10,14] = 275
$arr[10,15] = 4
$arr[11,0] = 3.7
$arr[11,1] = 6
$arr[11,2] = 305
$arr[11,3] = 5
$arr[11,4] = 4
$arr[11,5] = 7.8
$arr[11,6] = 212
$arr[11,7] = 378
$arr[11,8] = 1821
$arr[11,9] = 569
$arr[11,10] = 60
$arr[11,11] = 515
$arr[11,12] = 194
$arr[11,13] = 154
$arr[11,14] = 287
$arr[11,15] = 6
$arr[12,0] = 0
$arr[12,1] = 2
$arr[12,2] = 761
$arr[12,3] = 4
$arr[12,4] = 4
$arr[12,5] = 2.8
$arr[12,6] = 260
$arr[12,7] = 1050
$arr[12,8] = 2295
$arr[12,9] = 366
$arr[12,10] = 0
$arr[12,11] = 496
$arr[12,12] = 197
$arr[12,13] = 138
$arr[12,14] = 290
$arr[12,15] = 2
$arr[13,0] = 3
$arr[13,1] = 6
$arr[13,2] = 380
$arr[13,3] = 2
$arr[13,4] = 2
$arr[13,5] = 5.5
$arr[13,6] = 275
$arr[13,7] = 460
$arr[13,8] = 1584
$arr[13,9] = 310
$arr[13,10] = 70
$arr[13,11] = 448
$arr[13,12] = 192
$arr[13,13] = 131
$arr[13,14] = 262
$arr[13,15] = 6
$arr[14,0] = 2.4
$arr[14,1] = 6
$arr[14,2] = 195
$arr[14,3] = 2
$arr[14,4] = 2
$arr[14,5] = 7.1
$arr[14,6] = 238
$arr[14,7] = 225
$arr[14,8] = 1080
$arr[14,9] = 300
$arr[14,10] = 65
$arr[14,11] = 424
$arr[14,12] = 170
$arr[14,13] = 114
$arr[14,14] = 234
$arr[14,15] = 5
$arr[15,0] = 2.7
$arr[15,1] = 6
$arr[15,2] = 220
$arr[15,3] = 2
$arr[15,4] = 2
$arr[15,5] = 6.6
$arr[15,6] = 250
$arr[15,7] = 260
$arr[15,8] = 1260
$arr[15,9] = 260
$arr[15,10] = 64
$arr[15,11] = 432
$arr[15,12] = 178
$arr[15,13] = 129
$arr[15,14] = 242
$arr[15,15] = 5
$arr[16,0] = 2
$arr[16,1] = 4
$arr[16,2] = 160
$arr[16,3] = 4
$arr[16,4] = 2
$arr[16,5] = 6.8
$arr[16,6] = 221
$arr[16,7] = 221
$arr[16,8] = 900
$arr[16,9] = 198
$arr[16,10] = 68
$arr[16,11] = 419
$arr[16,12] = 186
$arr[16,13] = 111
$arr[16,14] = 244
$arr[16,15] = 5
$arr[17,0] = 5
$arr[17,1] = 8
$arr[17,2] = 306
$arr[17,3] = 4
$arr[17,4] = 2
$arr[17,5] = 6.3
$arr[17,6] = 250
$arr[17,7] = 460
$arr[17,8] = 1810
$arr[17,9] = 450
$arr[17,10] = 88
$arr[17,11] = 499
$arr[17,12] = 186
$arr[17,13] = 140
$arr[17,14] = 289
$arr[17,15] = 7
$arr[18,0] = 2
$arr[18,1] = 4
$arr[18,2] = 300
$arr[18,3] = 2
$arr[18,4] = 2
$arr[18,5] = 5.3
$arr[18,6] = 275
$arr[18,7] = 380
$arr[18,8] = 1335
$arr[18,9] = 150
$arr[18,10] = 54
$arr[18,11] = 438
$arr[18,12] = 180
$arr[18,13] = 128
$arr[18,14] = 248
$arr[18,15] = 6
$arr[19,0] = 2
$arr[19,1] = 4
$arr[19,2] = 210
$arr[19,3] = 5
$arr[19,4] = 5
$arr[19,5] = 6.9
$arr[19,6] = 240
$arr[19,7] = 280
$arr[19,8] = 1318
$arr[19,9] = 350
$arr[19,10] = 55
$arr[19,11] = 421
$arr[19,12] = 179
$arr[19,13] = 147
$arr[19,14] = 258
$arr[19,15] = 6
$arr[20,0] = 7
$arr[20,1] = 8
$arr[20,2] = 385
$arr[20,3] = 5
$arr[20,4] = 4
$arr[20,5] = 6.5
$arr[20,6] = 215
$arr[20,7] = 610
$arr[20,8] = 1740
$arr[20,9] = 813
$arr[20,10] = 91
$arr[20,11] = 542
$arr[20,12] = 203
$arr[20,13] = 138
$arr[20,14] = 302
$arr[20,15] = 4
$arr[21,0] = 6.5
$arr[21,1] = 12
$arr[21,2] = 819
$arr[21,3] = 2
$arr[21,4] = 2
$arr[21,5] = 2.8
$arr[21,6] = 355
$arr[21,7] = 720
$arr[21,8] = 1595
$arr[21,9] = 63
$arr[21,10] = 70
$arr[21,11] = 487
$arr[21,12] = 210
$arr[21,13] = 114
$arr[21,14] = 270
$arr[21,15] = 7
$arr[22,0] = 6.5
$arr[22,1] = 12
$arr[22,2] = 640
$arr[22,3] = 2
$arr[22,4] = 2
$arr[22,5] = 3.4
$arr[22,6] = 338
$arr[22,7] = 660
$arr[22,8] = 1830
$arr[22,9] = 140
$arr[22,10] = 100
$arr[22,11] = 461
$arr[22,12] = 206
$arr[22,13] = 114
$arr[22,14] = 267
$arr[22,15] = 6
$arr[23,0] = 5
$arr[23,1] = 8
$arr[23,2] = 600
$arr[23,3] = 4
$arr[23,4] = 4
$arr[23,5] = 3.7
$arr[23,6] = 322
$arr[23,7] = 700
$arr[23,8] = 1745
$arr[23,9] = 450
$arr[23,10] = 75
$arr[23,11] = 471
$arr[23,12] = 208
$arr[23,13] = 144
$arr[23,14] = 284
$arr[23,15] = 8
$arr[24,0] = 6.2
$arr[24,1] = 8
$arr[24,2] = 462
$arr[24,3] = 2
$arr[24,4] = 2
$arr[24,5] = 3.5
$arr[24,6] = 296
$arr[24,7] = 613
$arr[24,8] = 1730
$arr[24,9] = 357
$arr[24,10] = 70
$arr[24,11] = 463
$arr[24,12] = 193
$arr[24,13] = 124
$arr[24,14] = 272
$arr[24,15] = 8
$arr[25,0] = 3.5
$arr[25,1] = 6
$arr[25,2] = 656
$arr[25,3] = 2
$arr[25,4] = 2
$arr[25,5] = 3.2
$arr[25,6] = 347
$arr[25,7] = 746
$arr[25,8] = 1385
$arr[25,9] = 501
$arr[25,10] = 57
$arr[25,11] = 476
$arr[25,12] = 200
$arr[25,13] = 111
$arr[25,14] = 271
$arr[25,15] = 7
$arr[26,0] = 4.4
$arr[26,1] = 8
$arr[26,2] = 625
$arr[26,3] = 4
$arr[26,4] = 2
$arr[26,5] = 3.2
$arr[26,6] = 304
$arr[26,7] = 750
$arr[26,8] = 1900
$arr[26,9] = 420
$arr[26,10] = 68
$arr[26,11] = 487
$arr[26,12] = 191
$arr[26,13] = 214
$arr[26,14] = 283
$arr[26,15] = 8
$arr[27,0] = 6
$arr[27,1] = 12
$arr[27,2] = 600
$arr[27,3] = 3
$arr[27,4] = 2
$arr[27,5] = 3
$arr[27,6] = 362
$arr[27,7] = 627
$arr[27,8] = 1060
$arr[27,9] = 285
$arr[27,10] = 90
$arr[27,11] = 429
$arr[27,12] = 182
$arr[27,13] = 114
$arr[27,14] = 272
$arr[27,15] = 6
$arr[28,0] = 7
$arr[28,1] = 12
$arr[28,2] = 811
$arr[28,3] = 2
$arr[28,4] = 2
$arr[28,5] = 3
$arr[28,6] = 360
$arr[28,7] = 800
$arr[28,8] = 1350
$arr[28,9] = 0
$arr[28,10] = 120
$arr[28,11] = 481
$arr[28,12] = 206
$arr[28,13] = 124
$arr[28,14] = 280
$arr[28,15] = 6
$arr[29,0] = 3.5
$arr[29,1] = 6
$arr[29,2] = 573
$arr[29,3] = 2
$arr[29,4] = 2
$arr[29,5] = 3.8
$arr[29,6] = 307
$arr[29,7] = 550
$arr[29,8] = 1725
$arr[29,9] = 110
$arr[29,10] = 59
$arr[29,11] = 447
$arr[29,12] = 223
$arr[29,13] = 121
$arr[29,14] = 263
$arr[29,15] = 9
$arr[30,0] = 5.7
$arr[30,1] = 8
$arr[30,2] = 157
$arr[30,3] = 4
$arr[30,4] = 2
$arr[30,5] = 9.8
$arr[30,6] = 190
$arr[30,7] = 353
$arr[30,8] = 1626
$arr[30,9] = 312
$arr[30,10] = 79
$arr[30,11] = 500
$arr[30,12] = 186
$arr[30,13] = 125
$arr[30,14] = 175
$arr[30,15] = 3
$arr[31,0] = 3.8
$arr[31,1] = 6
$arr[31,2] = 700
$arr[31,3] = 2
$arr[31,4] = 2
$arr[31,5] = 2.8
$arr[31,6] = 340
$arr[31,7] = 700
$arr[31,8] = 1470
$arr[31,9] = 115
$arr[31,10] = 64
$arr[31,11] = 455
$arr[31,12] = 188
$arr[31,13] = 130
$arr[31,14] = 245
$arr[31,15] = 7
$arr[32,0] = 4.5
$arr[32,1] = 8
$arr[32,2] = 570
$arr[32,3] = 2
$arr[32,4] = 2
$arr[32,5] = 3.4
$arr[32,6] = 325
$arr[32,7] = 540
$arr[32,8] = 1380
$arr[32,9] = 229
$arr[32,10] = 86
$arr[32,11] = 453
$arr[32,12] = 194
$arr[32,13] = 121
$arr[32,14] = 265
$arr[32,15] = 7
$arr[33,0] = 5.2
$arr[33,1] = 10
$arr[33,2] = 570
$arr[33,3] = 2
$arr[33,4] = 2
$arr[33,5] = 3.4
$arr[33,6] = 325
$arr[33,7] = 540
$arr[33,8] = 1430
$arr[33,9] = 110
$arr[33,10] = 90
$arr[33,11] = 435
$arr[33,12] = 190
$arr[33,13] = 117
$arr[33,14] = 256
$arr[33,15] = 6
$arr[34,0] = 6.5
$arr[34,1] = 12
$arr[34,2] = 819
$arr[34,3] = 2
$arr[34,4] = 2
$arr[34,5] = 2.8
$arr[34,6] = 350
$arr[34,7] = 720
$arr[34,8] = 1601
$arr[34,9] = 140
$arr[34,10] = 85
$arr[34,11] = 498
$arr[34,12] = 210
$arr[34,13] = 113
$arr[34,14] = 270
$arr[34,15] = 7
$arr[35,0] = $null
$arr[35,1] = $null
$arr[35,2] = $null
$arr[35,3] = $null
$arr[35,4] = $null
$arr[35,5] = $null
$arr[35,6] = $null
$arr[35,7] = $null
$arr[35,8] = $null
$arr[35,9] = $null
$arr[35,10] = $null
$arr[35,11] = $null
$arr[35,12] = $null
$arr[35,13] = $null
$arr[35,14] = $null
$arr[35,15] = $null
$arr[36,0] = $null
$arr[36,1] = $null
$arr[36,2] = $null
$arr[36,3] = $null
$arr[36,4] = $null
$arr[36,5] = $null
$arr[36,6] = $null
$arr[36,7] = $null
$arr[36,8] = $null
$arr[36,9] = $null
$arr[36,10] = $null
$arr[36,11] = $null
$arr[36,12] = $null
$arr[36,13] = $null
$arr[36,14] = $null
$arr[36,15] = $null
$arr[37,0] = $null
$arr[37,1] = $null
$arr[37,2] = $null
$arr[37,3] = $null
$arr[37,4] = $null
$arr[37,5] = $null
$arr[37,6] = $null
$arr[37,7] = $null
$arr[37,8] = $null
$arr[37,9] = $null
$arr[37,10] = $null
$arr[37,11] = $null
$arr[37,12] = $null
$arr[37,13] = $null
$arr[37,14] = $null
$arr[37,15] = $null
$arr[38,0] = $null
$arr[38,1] = $null
$arr[38,2] = $null
$arr[38,3] = $null
$arr[38,4] = $null
$arr[38,5] = $null
$arr[38,6] = $null
$arr[38,7] = $null
$arr[38,8] = $null
$arr[38,9] = $null
$arr[38,10] = $null
$arr[38,11] = $null
$arr[38,12] = $null
$arr[38,13] = $null
$arr[38,14] = $null
$arr[38,15] = $null
$arr[39,0] = $null
$arr[39,1] = $null
$arr[39,2] = $null
$arr[39,3] = $null
$arr[39,4] = $null
$arr[39,5] = $null
$arr[39,6] = $null
$arr[39,7] = $null
$arr[39,8] = $null
$arr[39,9] = $null
$arr[39,10] = $null
$arr[39,11] = $null
$arr[39,12] = $null
$arr[39,13] = $null
$arr[39,14] = $null
$arr[39,15] = $null
$arr[40,0] = $null
$arr[40,1] = $null
$arr[40,2] = $null
$arr[40,3] = $null
$arr[40,4] = $null
$arr[40,5] = $null
$arr[40,6] = $null
$arr[40,7] = $null
$arr[40,8] = $null
$arr[40,9] = $null
$arr[40,10] = $null
$arr[40,11] = $null
$arr[40,12] = $null
$arr[40,13] = $null
$arr[40,14] = $null
$arr[40,15] = $null
$arr[41,0] = $null
$arr[41,1] = $null
$arr[41,2] = $null
$arr[41,3] = $null
$arr[41,4] = $null
$arr[41,5] = $null
$arr[41,6] = $null
$arr[41,7] = $null
$arr[41,8] = $null
$arr[41,9] = $null
$arr[41,10] = $null
$arr[41,11] = $null
$arr[41,12] = $null
$arr[41,13] = $null
$arr[41,14] = $null
$arr[41,15] = $null

$ws.Range("E167:T208").Value = $arr

# Update the sheet view: scroll position reset, new active selection
[void]$ws.Activate()
[void]$ws.Range("M15").Select()
